# Update the "events" sheet so that regular visit events are primarily
# identified by an explicit event_id (VIS1..VIS6) instead of relying on a
# single event_id_pattern ("^VIS[[:digit:]]+$") / ("^FU") row.
#
# Columns on the "events" table (Table913): A=event_id, B=event_id_pattern,
# C=is_regular_visit, D=event_label_custom, E=event_name_custom,
# F=is_baseline_event

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("events")

# --- Clear the previous data rows (2:7) so no stray cells are left over ---
$ws.Range("A2:F7").ClearContents()

# --- Row 2: SCR (unchanged) ---
$ws.Range("A2").Value = "SCR"
$ws.Range("C2").Value = $true
$ws.Range("E2").Value = "Screening"
$ws.Range("F2").Value = $true

# --- Rows 3-8: explicit VIS1..VIS6 events replacing the old pattern-based
#     "^VIS[[:digit:]]+$" row ---
$visitRows = 3..8
$visitNum = 1
foreach ($r in $visitRows) {
    $ws.Range("A$r").Value = "VIS$visitNum"
    $ws.Range("C$r").Value = $true
    $ws.Range("D$r").Value = "V$visitNum"
    $ws.Range("E$r").Value = "Visit $visitNum"
    $ws.Range("F$r").Value = $false
    $visitNum++
}

# --- Row 9: EOT (unchanged values, shifted down from old row 4) ---
$ws.Range("A9").Value = "EOT"
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = "EoT"
$ws.Range("E9").Value = "EoT"
$ws.Range("F9").Value = $false

# --- Row 10: EXIT (unchanged values, shifted down from old row 6) ---
$ws.Range("A10").Value = "EXIT"
$ws.Range("C10").Value = $false
$ws.Range("E10").Value = "Exit"
$ws.Range("F10").Value = $false

# --- Row 11: Unscheduled visit, still identified via event_id_pattern
#     (shifted down from old row 7) ---
$ws.Range("B11").Value = "^UN"
$ws.Range("C11").Value = $false
$ws.Range("E11").Value = "Unscheduled visit"
$ws.Range("F11").Value = $false

# --- Resize the table / list object to cover the new data extent ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F11"))

# --- Column width tweaks to better fit the new (shorter) event_id values
#     and slightly wider event_id_pattern column ---
$ws.Columns.Item(1).ColumnWidth = 14.8
$ws.Columns.Item(2).ColumnWidth = 17.25

# --- Leave the cursor where the author's session ended up ---
$ws.Range("D9").Select()
